$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus/vm_pu.xlsx values for the 380 kV case (commit: "case with 380 kV done")
# Column B (slack bus voltage setpoint) changes from 1.05 -> 1.02 p.u. for every row,
# and the resulting bus voltage magnitudes (columns C-F, I-N) are updated to the
# newly recomputed power-flow results for rows 2-25.

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031539701158009
$ws.Range("D2").Value = 1.040615143088023
$ws.Range("E2").Value = 1.04951739431286
$ws.Range("F2").Value = 1.054037019861213
$ws.Range("I2").Value = 1.036889344208608
$ws.Range("J2").Value = 1.036674867731149
$ws.Range("K2").Value = 1.043396989555646
$ws.Range("L2").Value = 1.052274231972689
$ws.Range("M2").Value = 1.056781335913599
$ws.Range("N2").Value = 1.016222365518382

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032276224986706
$ws.Range("D3").Value = 1.041193382845385
$ws.Range("E3").Value = 1.050351816297958
$ws.Range("F3").Value = 1.054827055175178
$ws.Range("I3").Value = 1.037027875723292
$ws.Range("J3").Value = 1.037054393523527
$ws.Range("K3").Value = 1.043786299837397
$ws.Range("L3").Value = 1.052920858718145
$ws.Range("M3").Value = 1.057384594018029
$ws.Range("N3").Value = 1.016348658099686

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032753464619126
$ws.Range("D4").Value = 1.041568156206712
$ws.Range("E4").Value = 1.050893256486715
$ws.Range("F4").Value = 1.055339496222175
$ws.Range("I4").Value = 1.037116649484956
$ws.Range("J4").Value = 1.037299929080495
$ws.Range("K4").Value = 1.044038120685944
$ws.Range("L4").Value = 1.053340132074985
$ws.Range("M4").Value = 1.057775517428688
$ws.Range("N4").Value = 1.016430344813956

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032954252080042
$ws.Range("D5").Value = 1.041725855881217
$ws.Range("E5").Value = 1.051121238082271
$ws.Range("F5").Value = 1.055555220002444
$ws.Range("I5").Value = 1.037153762062105
$ws.Range("J5").Value = 1.037403140642303
$ws.Range("K5").Value = 1.044143963499998
$ws.Range("L5").Value = 1.053516599303413
$ws.Range("M5").Value = 1.057939997506415
$ws.Range("N5").Value = 1.01646467754376

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032987974264524
$ws.Range("D6").Value = 1.041752342794231
$ws.Range("E6").Value = 1.051159538251231
$ws.Range("F6").Value = 1.055591458125544
$ws.Range("I6").Value = 1.037159981213674
$ws.Range("J6").Value = 1.03742046960074
$ws.Range("K6").Value = 1.044161733604676
$ws.Range("L6").Value = 1.053546240879403
$ws.Range("M6").Value = 1.057967622350726
$ws.Range("N6").Value = 1.016470441656341

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032756146940803
$ws.Range("D7").Value = 1.041570262830353
$ws.Range("E7").Value = 1.050896301375567
$ws.Range("F7").Value = 1.055342377583127
$ws.Range("I7").Value = 1.037117146202898
$ws.Range("J7").Value = 1.037301308245581
$ws.Range("K7").Value = 1.044039535053618
$ws.Range("L7").Value = 1.053342489236626
$ws.Range("M7").Value = 1.057777714687898
$ws.Range("N7").Value = 1.016430803603015

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03178847511018
$ws.Range("D8").Value = 1.040810433713285
$ws.Range("E8").Value = 1.04979907647591
$ws.Range("F8").Value = 1.054303758888585
$ws.Range("I8").Value = 1.0369363403249
$ws.Range("J8").Value = 1.03680313837947
$ws.Range("K8").Value = 1.043528576373884
$ws.Range("L8").Value = 1.052492582815253
$ws.Range("M8").Value = 1.056985089659043
$ws.Range("N8").Value = 1.016265053251999

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030088451109291
$ws.Range("D9").Value = 1.03947630259239
$ws.Range("E9").Value = 1.0478773123793
$ws.Range("F9").Value = 1.052483138266661
$ws.Range("I9").Value = 1.036611142012435
$ws.Range("J9").Value = 1.035925027096857
$ws.Range("K9").Value = 1.042627581770083
$ws.Range("L9").Value = 1.051001625490694
$ws.Range("M9").Value = 1.055592867128214
$ws.Range("N9").Value = 1.015972747196394

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028958673098743
$ws.Range("D10").Value = 1.038590214852859
$ws.Range("E10").Value = 1.046604123254159
$ws.Range("F10").Value = 1.051275948597394
$ws.Range("I10").Value = 1.036389953015349
$ws.Range("J10").Value = 1.035339513689774
$ws.Range("K10").Value = 1.04202658435037
$ws.Range("L10").Value = 1.050012257516457
$ws.Range("M10").Value = 1.054667841394622
$ws.Range("N10").Value = 1.015777747347153

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028470338658366
$ws.Range("D11").Value = 1.038207344243892
$ws.Range("E11").Value = 1.04605474025041
$ws.Range("F11").Value = 1.050754805296731
$ws.Range("I11").Value = 1.036293143816779
$ws.Range("J11").Value = 1.035085971177722
$ws.Range("K11").Value = 1.041766284049365
$ws.Range("L11").Value = 1.049584965107528
$ws.Range("M11").Value = 1.054268058569425
$ws.Range("N11").Value = 1.015693285426056

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028289081762494
$ws.Range("D12").Value = 1.038065252940582
$ws.Range("E12").Value = 1.045850965097159
$ws.Range("F12").Value = 1.050561468790885
$ws.Range("I12").Value = 1.036257030135143
$ws.Range("J12").Value = 1.034991793822802
$ws.Range("K12").Value = 1.041669588665438
$ws.Range("L12").Value = 1.049426418387104
$ws.Range("M12").Value = 1.054119677387066
$ws.Range("N12").Value = 1.015661909116881

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028327955964581
$ws.Range("D13").Value = 1.038095726369048
$ws.Range("E13").Value = 1.045894662392912
$ws.Range("F13").Value = 1.05060292926145
$ws.Range("I13").Value = 1.036264783618938
$ws.Range("J13").Value = 1.035011995180883
$ws.Range("K13").Value = 1.041690330504611
$ws.Range("L13").Value = 1.049460419532498
$ws.Range("M13").Value = 1.054151500377041
$ws.Range("N13").Value = 1.015668639587618

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028455353197108
$ws.Range("D14").Value = 1.038195596393543
$ws.Range("E14").Value = 1.04603789020055
$ws.Range("F14").Value = 1.050738819145643
$ws.Range("I14").Value = 1.036290161794687
$ws.Range("J14").Value = 1.035078186448203
$ws.Range("K14").Value = 1.041758291339098
$ws.Range("L14").Value = 1.049571856133792
$ws.Range("M14").Value = 1.054255790958288
$ws.Range("N14").Value = 1.015690691917011

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028533864422585
$ws.Range("D15").Value = 1.038257146088204
$ws.Range("E15").Value = 1.046126176095819
$ws.Range("F15").Value = 1.050822577168447
$ws.Range("I15").Value = 1.036305777684609
$ws.Range("J15").Value = 1.035118969062319
$ws.Range("K15").Value = 1.041800163199873
$ws.Range("L15").Value = 1.04964053831283
$ws.Range("M15").Value = 1.054320063248907
$ws.Range("N15").Value = 1.015704278649676

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028991100683898
$ws.Range("D16").Value = 1.038615641977566
$ws.Range("E16").Value = 1.046640624561265
$ws.Range("F16").Value = 1.051310568598365
$ws.Range("I16").Value = 1.036396356224688
$ws.Range("J16").Value = 1.0353563403456
$ws.Range("K16").Value = 1.042043858387449
$ws.Range("L16").Value = 1.05004063906219
$ws.Range("M16").Value = 1.054694389823408
$ws.Range("N16").Value = 1.015783352305364

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029278146426498
$ws.Range("D17").Value = 1.038840735586996
$ws.Range("E17").Value = 1.046963839129033
$ws.Range("F17").Value = 1.051617096755615
$ws.Range("I17").Value = 1.036452897739311
$ws.Range("J17").Value = 1.035505234980642
$ws.Range("K17").Value = 1.042196705687958
$ws.Range("L17").Value = 1.050291910313422
$ws.Range("M17").Value = 1.05492939962934
$ws.Range("N17").Value = 1.015832946571988

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029445658859049
$ws.Range("D18").Value = 1.03897210701838
$ws.Range("E18").Value = 1.047152549460981
$ws.Range("F18").Value = 1.051796041525309
$ws.Range("I18").Value = 1.036485777737819
$ws.Range("J18").Value = 1.035592081518768
$ws.Range("K18").Value = 1.042285852618417
$ws.Range("L18").Value = 1.050438579600119
$ws.Range("M18").Value = 1.055066550021248
$ws.Range("N18").Value = 1.015861871557515

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029502790372565
$ws.Range("D19").Value = 1.03901691447377
$ws.Range("E19").Value = 1.047216926074463
$ws.Range("F19").Value = 1.051857082793893
$ws.Range("I19").Value = 1.036496972032295
$ws.Range("J19").Value = 1.035621693702839
$ws.Range("K19").Value = 1.04231624830374
$ws.Range("L19").Value = 1.050488608131583
$ws.Range("M19").Value = 1.05511332711106
$ws.Range("N19").Value = 1.01587173378985

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029247340482796
$ws.Range("D20").Value = 1.03881657707637
$ws.Range("E20").Value = 1.046929142139706
$ws.Range("F20").Value = 1.051584193438387
$ws.Range("I20").Value = 1.036446841678025
$ws.Range("J20").Value = 1.035489260108264
$ws.Range("K20").Value = 1.042180307261986
$ws.Range("L20").Value = 1.050264940203403
$ws.Range("M20").Value = 1.054904177721777
$ws.Range("N20").Value = 1.015827625832883

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028417834220238
$ws.Range("D21").Value = 1.03816618372448
$ws.Range("E21").Value = 1.045995705153239
$ws.Range("F21").Value = 1.050698796337054
$ws.Range("I21").Value = 1.036282692811614
$ws.Range("J21").Value = 1.035058694769503
$ws.Range("K21").Value = 1.041738278783893
$ws.Range("L21").Value = 1.049539036160064
$ws.Range("M21").Value = 1.054225076765472
$ws.Range("N21").Value = 1.015684198145481

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027897056915362
$ws.Range("D22").Value = 1.037757973381043
$ws.Range("E22").Value = 1.045410496138455
$ws.Range("F22").Value = 1.050143497794272
$ws.Range("I22").Value = 1.036178592762639
$ws.Range("J22").Value = 1.034787979727302
$ws.Range("K22").Value = 1.04146031083541
$ws.Range("L22").Value = 1.049083608481197
$ws.Range("M22").Value = 1.053798770867744
$ws.Range("N22").Value = 1.015594000078662

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028173057516819
$ws.Range("D23").Value = 1.03797430470112
$ws.Range("E23").Value = 1.045720566525004
$ws.Range("F23").Value = 1.050437739945609
$ws.Range("I23").Value = 1.036233862577888
$ws.Range("J23").Value = 1.034931490594833
$ws.Range("K23").Value = 1.041607670894565
$ws.Range("L23").Value = 1.049324946076896
$ws.Range("M23").Value = 1.054024699386208
$ws.Range("N23").Value = 1.015641817461363

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029261260104976
$ws.Range("D24").Value = 1.038827493026275
$ws.Range("E24").Value = 1.046944819645152
$ws.Range("F24").Value = 1.051599060561539
$ws.Range("I24").Value = 1.036449578459629
$ws.Range("J24").Value = 1.035496478469297
$ws.Range("K24").Value = 1.042187717024522
$ws.Range("L24").Value = 1.050277126504795
$ws.Range("M24").Value = 1.054915574190378
$ws.Range("N24").Value = 1.015830030053679

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030527327179484
$ws.Range("D25").Value = 1.039820628476757
$ws.Range("E25").Value = 1.048372736204997
$ws.Range("F25").Value = 1.052952666137451
$ws.Range("I25").Value = 1.036695990598712
$ws.Range("J25").Value = 1.03615206422587
$ws.Range("K25").Value = 1.042860575447767
$ws.Range("L25").Value = 1.051386270195105
$ws.Range("M25").Value = 1.055952247532487
$ws.Range("N25").Value = 1.016048340081608
